$p = $ppt.ActivePresentation
$s = $p.Slides.Add(3, 2)
$s.Shapes.Item(1).Delete()
$shape = $s.Shapes.Item(1)
$tr = $shape.TextFrame.TextRange
$tr.Text = "`rThank you"

$run2 = $tr.Characters(2, 9)
Write-Output ("run2 text: [" + $run2.Text + "]")
$run2.Font.NameFarEast = "BM KIRANGHAERANG OTF"

Write-Output "done"
